$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (B1/C1) reuse the same style as A1 (bold title style) ---
$ws.Range("A1").Copy($ws.Range("B1"))
$ws.Range("A1").Copy($ws.Range("C1"))
$ws.Range("B1").Value = "Latitude"
$ws.Range("C1").Value = "Longitude"

# --- Column widths for B and C (target stored width = 30) ---
$ws.Columns.Item(2).ColumnWidth = 29.15
$ws.Columns.Item(3).ColumnWidth = 29.15

# --- Latitude / Longitude data cells, left-aligned ---
$ws.Range("B2").Value = -9.549443999999999
$ws.Range("C2").Value = -35.72652
$ws.Range("B2:C2").HorizontalAlignment = -4131
$ws.Range("B3").Value = -9.5476425
$ws.Range("C3").Value = -35.7293331
$ws.Range("B3:C3").HorizontalAlignment = -4131
$ws.Range("B4").Value = -9.549807899999999
$ws.Range("C4").Value = -35.7331086
$ws.Range("B4:C4").HorizontalAlignment = -4131
$ws.Range("B5").Value = -9.5520394
$ws.Range("C5").Value = -35.7353938
$ws.Range("B5:C5").HorizontalAlignment = -4131
$ws.Range("B6").Value = -9.160664499999999
$ws.Range("C6").Value = -35.5350591
$ws.Range("B6:C6").HorizontalAlignment = -4131
$ws.Range("B7").Value = -9.556360699999999
$ws.Range("C7").Value = -35.7399394
$ws.Range("B7:C7").HorizontalAlignment = -4131
$ws.Range("B8").Value = -9.5580718
$ws.Range("C8").Value = -35.74165929999999
$ws.Range("B8:C8").HorizontalAlignment = -4131
$ws.Range("B9").Value = -9.559123999999999
$ws.Range("C9").Value = -35.7427718
$ws.Range("B9:C9").HorizontalAlignment = -4131
$ws.Range("B10").Value = -9.560942499999999
$ws.Range("C10").Value = -35.7448797
$ws.Range("B10:C10").HorizontalAlignment = -4131
$ws.Range("B11").Value = -9.5613946
$ws.Range("C11").Value = -35.74615989999999
$ws.Range("B11:C11").HorizontalAlignment = -4131
$ws.Range("B12").Value = -9.565182399999999
$ws.Range("C12").Value = -35.777396
$ws.Range("B12:C12").HorizontalAlignment = -4131
$ws.Range("B13").Value = -9.546670300000001
$ws.Range("C13").Value = -35.7462117
$ws.Range("B13:C13").HorizontalAlignment = -4131
$ws.Range("B14").Value = -9.547639199999999
$ws.Range("C14").Value = -35.7478589
$ws.Range("B14:C14").HorizontalAlignment = -4131
$ws.Range("B15").Value = -9.5487631
$ws.Range("C15").Value = -35.7502234
$ws.Range("B15:C15").HorizontalAlignment = -4131
$ws.Range("B16").Value = -9.549141199999999
$ws.Range("C16").Value = -35.7518819
$ws.Range("B16:C16").HorizontalAlignment = -4131
$ws.Range("B17").Value = -9.5485849
$ws.Range("C17").Value = -35.75471599999999
$ws.Range("B17:C17").HorizontalAlignment = -4131
$ws.Range("B18").Value = -9.547731799999999
$ws.Range("C18").Value = -35.7532515
$ws.Range("B18:C18").HorizontalAlignment = -4131
$ws.Range("B19").Value = -9.546125
$ws.Range("C19").Value = -35.7526961
$ws.Range("B19:C19").HorizontalAlignment = -4131
$ws.Range("B21").Value = -9.542885099999999
$ws.Range("C21").Value = -35.7528852
$ws.Range("B21:C21").HorizontalAlignment = -4131
$ws.Range("B22").Value = -9.547007499999999
$ws.Range("C22").Value = -35.7520547
$ws.Range("B22:C22").HorizontalAlignment = -4131
$ws.Range("B23").Value = -9.5392022
$ws.Range("C23").Value = -35.7541773
$ws.Range("B23:C23").HorizontalAlignment = -4131
$ws.Range("B24").Value = -9.539851499999999
$ws.Range("C24").Value = -35.7560938
$ws.Range("B24:C24").HorizontalAlignment = -4131
$ws.Range("B25").Value = -9.5399694
$ws.Range("C25").Value = -35.7563304
$ws.Range("B25:C25").HorizontalAlignment = -4131
$ws.Range("B28").Value = -9.5416533
$ws.Range("C28").Value = -35.7595592
$ws.Range("B28:C28").HorizontalAlignment = -4131
$ws.Range("B29").Value = 18.3654432
$ws.Range("C29").Value = -66.13259719999999
$ws.Range("B29:C29").HorizontalAlignment = -4131
$ws.Range("B30").Value = -9.5505864
$ws.Range("C30").Value = -35.7593578
$ws.Range("B30:C30").HorizontalAlignment = -4131
$ws.Range("B31").Value = -9.5516975
$ws.Range("C31").Value = -35.7612876
$ws.Range("B31:C31").HorizontalAlignment = -4131
$ws.Range("B32").Value = -9.5474643
$ws.Range("C32").Value = -35.754939
$ws.Range("B32:C32").HorizontalAlignment = -4131
$ws.Range("B33").Value = -9.5470734
$ws.Range("C33").Value = -35.7542347
$ws.Range("B33:C33").HorizontalAlignment = -4131
$ws.Range("B34").Value = -9.5500852
$ws.Range("C34").Value = -35.7592794
$ws.Range("B34:C34").HorizontalAlignment = -4131
$ws.Range("B35").Value = -9.5435246
$ws.Range("C35").Value = -35.7627701
$ws.Range("B35:C35").HorizontalAlignment = -4131
$ws.Range("B36").Value = -9.552133999999999
$ws.Range("C36").Value = -35.77124999999999
$ws.Range("B36:C36").HorizontalAlignment = -4131
$ws.Range("B37").Value = -9.552133999999999
$ws.Range("C37").Value = -35.77124999999999
$ws.Range("B37:C37").HorizontalAlignment = -4131
$ws.Range("B38").Value = -9.552133999999999
$ws.Range("C38").Value = -35.77124999999999
$ws.Range("B38:C38").HorizontalAlignment = -4131
$ws.Range("B39").Value = -9.555642599999999
$ws.Range("C39").Value = -35.7764813
$ws.Range("B39:C39").HorizontalAlignment = -4131
$ws.Range("B40").Value = -9.552133999999999
$ws.Range("C40").Value = -35.77124999999999
$ws.Range("B40:C40").HorizontalAlignment = -4131
$ws.Range("B41").Value = -9.5573503
$ws.Range("C41").Value = -35.7825882
$ws.Range("B41:C41").HorizontalAlignment = -4131
$ws.Range("B42").Value = -9.5564714
$ws.Range("C42").Value = -35.7812506
$ws.Range("B42:C42").HorizontalAlignment = -4131
$ws.Range("B43").Value = -9.5526211
$ws.Range("C43").Value = -35.790461
$ws.Range("B43:C43").HorizontalAlignment = -4131
$ws.Range("B44").Value = -9.5588791
$ws.Range("C44").Value = -35.7853729
$ws.Range("B44:C44").HorizontalAlignment = -4131
$ws.Range("B45").Value = -9.556795299999999
$ws.Range("C45").Value = -35.7865437
$ws.Range("B45:C45").HorizontalAlignment = -4131
$ws.Range("B46").Value = -9.5581777
$ws.Range("C46").Value = -35.78896170000001
$ws.Range("B46:C46").HorizontalAlignment = -4131
$ws.Range("B47").Value = -9.555773799999999
$ws.Range("C47").Value = -35.7904942
$ws.Range("B47:C47").HorizontalAlignment = -4131
$ws.Range("B48").Value = -9.553751999999999
$ws.Range("C48").Value = -35.7909765
$ws.Range("B48:C48").HorizontalAlignment = -4131
$ws.Range("B49").Value = -9.551980499999999
$ws.Range("C49").Value = -35.7901472
$ws.Range("B49:C49").HorizontalAlignment = -4131
$ws.Range("B50").Value = -9.551114
$ws.Range("C50").Value = -35.7873119
$ws.Range("B50:C50").HorizontalAlignment = -4131
$ws.Range("B51").Value = -9.546662
$ws.Range("C51").Value = -35.791912
$ws.Range("B51:C51").HorizontalAlignment = -4131
$ws.Range("B52").Value = -9.5507568
$ws.Range("C52").Value = -35.786958
$ws.Range("B52:C52").HorizontalAlignment = -4131
$ws.Range("B53").Value = -22.9111438
$ws.Range("C53").Value = -43.1648755
$ws.Range("B53:C53").HorizontalAlignment = -4131
$ws.Range("B55").Value = -9.5377907
$ws.Range("C55").Value = -35.81515340000001
$ws.Range("B55:C55").HorizontalAlignment = -4131
$ws.Range("B56").Value = -9.540084199999999
$ws.Range("C56").Value = -35.7965546
$ws.Range("B56:C56").HorizontalAlignment = -4131
$ws.Range("B58").Value = -9.543493399999999
$ws.Range("C58").Value = -35.7870338
$ws.Range("B58:C58").HorizontalAlignment = -4131
$ws.Range("B59").Value = -9.542950099999999
$ws.Range("C59").Value = -35.7860535
$ws.Range("B59:C59").HorizontalAlignment = -4131
$ws.Range("B60").Value = -9.5418305
$ws.Range("C60").Value = -35.7832275
$ws.Range("B60:C60").HorizontalAlignment = -4131

Write-Host "done"
